$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-CellText($row, $col, $oldText, $newText) {
    $cellRange = $t.Cell($row, $col).Range
    # Replace = 1 (wdReplaceOne) so the substitution stays confined to this cell's range
    # even though the same text occurs in other cells of the table.
    $cellRange.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 1)
}

Set-CellText 1 1 "80÷7=11, 3"  "98÷8=12, 2"
Set-CellText 1 2 "45÷4=11, 1"  "56÷8=7, 0"
Set-CellText 1 3 "62÷8=7, 6"   "30÷3=10, 0"
Set-CellText 1 4 "47÷6=7, 5"   "44÷2=22, 0"
Set-CellText 1 5 "89÷6=14, 5"  "37÷7=5, 2"

Set-CellText 5 1 "24÷5=4, 4"   "31÷6=5, 1"
Set-CellText 5 2 "65÷3=21, 2"  "97÷3=32, 1"
Set-CellText 5 3 "29÷6=4, 5"   "23÷3=7, 2"
Set-CellText 5 4 "59÷6=9, 5"   "87÷5=17, 2"
Set-CellText 5 5 "64÷7=9, 1"   "32÷3=10, 2"

Set-CellText 9 1 "29÷6=4, 5"   "43÷5=8, 3"
Set-CellText 9 2 "30÷4=7, 2"   "93÷3=31, 0"
Set-CellText 9 3 "16÷7=2, 2"   "95÷2=47, 1"
Set-CellText 9 4 "33÷7=4, 5"   "43÷8=5, 3"
Set-CellText 9 5 "33÷6=5, 3"   "12÷8=1, 4"

Set-CellText 13 1 "14÷5=2, 4"  "55÷7=7, 6"
Set-CellText 13 2 "25÷3=8, 1"  "51÷5=10, 1"
Set-CellText 13 3 "39÷8=4, 7"  "86÷4=21, 2"
Set-CellText 13 4 "72÷3=24, 0" "77÷5=15, 2"
Set-CellText 13 5 "10÷6=1, 4"  "77÷8=9, 5"

Set-CellText 17 1 "40÷4=10, 0" "77÷7=11, 0"
Set-CellText 17 2 "15÷3=5, 0"  "86÷5=17, 1"
Set-CellText 17 3 "82÷3=27, 1" "84÷6=14, 0"
Set-CellText 17 4 "94÷6=15, 4" "21÷4=5, 1"
Set-CellText 17 5 "36÷7=5, 1"  "64÷9=7, 1"
